# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Updates the Serbia MSME country-indicator figures to more precise values:
#   Enterprises density (per 1000 people): 10.1 -> 10.06, 1.5 -> 1.53, 11.6 -> 11.58
#   Employment (% of total):               39.5 -> 39.46, 60 -> 59.96
#   Enterprises (% of total):              86.3 -> 86.31, 13.1 -> 13.11, 99.4 -> 99.42

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep a reference cell that already carries the plain/default style (s="0")
# used by every data cell in these rows, so we can restore it below after
# Excel auto-applies a "quote prefix" style to the newly typed numeric text.
$refCell = $ws.Range("B12")

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    # Leading apostrophe forces Excel to store the value as text (keeping it
    # a shared string) instead of silently coercing the numeric-looking
    # string into a numeric cell.
    $cell.Value = "'" + $Text
    # Re-apply the original (default) cell formatting so the style index
    # stays the same as before the edit (undoing the automatic "text
    # quote prefix" style Excel assigns to the cell above).
    $refCell.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# Row 13: Enterprises density (per 1000 people)
Set-TextValue "B13" "10.06"
Set-TextValue "C13" "1.53"
Set-TextValue "D13" "11.58"

# Row 14: Employment (% of total)
Set-TextValue "C14" "39.46"
Set-TextValue "D14" "59.96"

# Row 16: Enterprises (% of total)
Set-TextValue "B16" "86.31"
Set-TextValue "C16" "13.11"
Set-TextValue "D16" "99.42"

$excel.CutCopyMode = 0
